# Update CVD (Manufacturing Voluntary Turnover) ytd figures from 0.0639 to
# 0.0776 across all location sheets, plus the related Lavergne Tennessee
# (Internal Fill Rate) and Pharr Texas / Indianapolis Indiana monthly-detail
# corrections that came along with the refreshed CVD export.

$wb = $excel.ActiveWorkbook

# --- Chino California ---------------------------------------------------
$ws = $wb.Worksheets.Item("Chino California")
$ws.Range("E5").Value = 0.0776
$ws.Range("E6").Value = 0.0776
$ws.Range("E7").Value = 0.0776

# --- Mississauga Canada ---------------------------------------------------
$ws = $wb.Worksheets.Item("Mississauga Canada")
$ws.Range("E7").Value = 0.0776
$ws.Range("E8").Value = 0.0776
$ws.Range("E9").Value = 0.0776

# --- Edmonton EDM Canada ---------------------------------------------------
$ws = $wb.Worksheets.Item("Edmonton EDM Canada")
$ws.Range("E3").Value = 0.0776
$ws.Range("E4").Value = 0.0776
$ws.Range("E5").Value = 0.0776

# --- Mississauga Mid-Way Canada --------------------------------------------
$ws = $wb.Worksheets.Item("Mississauga Mid-Way Canada")
$ws.Range("E2").Value = 0.0776
$ws.Range("E3").Value = 0.0776
$ws.Range("E4").Value = 0.0776

# --- Monterrey Rbm Mexico ---------------------------------------------------
$ws = $wb.Worksheets.Item("Monterrey Rbm Mexico")
$ws.Range("E2").Value = 0.0776
$ws.Range("E3").Value = 0.0776

# --- Montreal Canada ---------------------------------------------------
$ws = $wb.Worksheets.Item("Montreal Canada")
$ws.Range("E2").Value = 0.0776
$ws.Range("E3").Value = 0.0776

# --- El Paso Texas - EPC ---------------------------------------------------
$ws = $wb.Worksheets.Item("El Paso Texas - EPC")
$ws.Range("L7").ClearContents()
$ws.Range("E8").Value = 0.0776
$ws.Range("E9").Value = 0.0776
$ws.Range("E10").Value = 0.0776

# --- Florence Kentucky ---------------------------------------------------
$ws = $wb.Worksheets.Item("Florence Kentucky")
$ws.Range("E7").Value = 0.0776
$ws.Range("E8").Value = 0.0776
$ws.Range("E9").Value = 0.0776

# --- Indianapolis Indiana ---------------------------------------------------
$ws = $wb.Worksheets.Item("Indianapolis Indiana")
$ws.Range("L7").ClearContents()
$ws.Range("E8").Value = 0.0776
$ws.Range("E9").Value = 0.0776
$ws.Range("E10").Value = 0.0776
$ws.Range("L10").Value = 0.0769

# --- Lavergne Tennessee ---------------------------------------------------
$ws = $wb.Worksheets.Item("Lavergne Tennessee")
$ws.Range("E5").Value = 0.666666666666667
$ws.Range("E6").Value = 0.666666666666667
$ws.Range("E7").Value = 0.666666666666667
$ws.Range("L7").Value = 1
$ws.Range("M7:W7").Value = 0.666666666666667
$ws.Range("E8").Value = 0.0776
$ws.Range("E9").Value = 0.0776
$ws.Range("E10").Value = 0.0776

# --- Pharr Texas ---------------------------------------------------
$ws = $wb.Worksheets.Item("Pharr Texas")
$ws.Range("L4").Value = 0.1429
$ws.Range("E6").Value = 0.0776
$ws.Range("E7").Value = 0.0776
$ws.Range("E8").Value = 0.0776
$ws.Range("L8:W8").Value = 0
